$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "34.854.62"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "1.804.32"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  +0.14%  "
$c = $ws.Range("D5")
$c.Value = "'231.51"
$c.ClearFormats()
$c = $ws.Range("D6")
$c.Value = "'0.602"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -8.91%  "
$ws.Range("E9").Value = "  +3.33%  "
$c = $ws.Range("D10")
$c.Value = "'0.0677"
$c.ClearFormats()
$ws.Range("E10").Value = "  -2.77%  "
$c = $ws.Range("D11")
$c.Value = "'0.0990"
$c.ClearFormats()
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").Value = "2.063.90"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").Value = "1.801.00"
$ws.Range("E13").Value = "  -3.28%  "
$c = $ws.Range("D14")
$c.Value = "'0.655"
$c.ClearFormats()
$ws.Range("E14").Value = "  -3.85%  "
$c = $ws.Range("D15")
$c.Value = "'10.81"
$c.ClearFormats()
$ws.Range("E15").Value = "  -6.34%  "
$c = $ws.Range("D16")
$c.Value = "'4.54"
$c.ClearFormats()
$ws.Range("E16").Value = "  -5.16%  "
$ws.Range("D17").Value = "34.770.96"
$ws.Range("E17").Value = "  -2.54%  "
$c = $ws.Range("D18")
$c.Value = "'68.92"
$c.ClearFormats()
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("E19").Value = "  -2.91%  "
$c = $ws.Range("D20")
$c.Value = "'238.46"
$c.ClearFormats()
$ws.Range("E20").Value = "  -3.07%  "
$c = $ws.Range("D21")
$c.Value = "'11.69"
$c.ClearFormats()
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("E23").Value = "  +0.17%  "
$c = $ws.Range("D24")
$c.Value = "'2.21"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.76%  "
$c = $ws.Range("D25")
$c.Value = "'171.15"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.44%  "
$c = $ws.Range("D26")
$c.Value = "'7.71"
$c.ClearFormats()
$ws.Range("E26").Value = "  -3.40%  "
$c = $ws.Range("D27")
$c.Value = "'17.21"
$c.ClearFormats()
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("E28").Value = "  -3.76%  "
$c = $ws.Range("D29")
$c.Value = "'1.54"
$c.ClearFormats()
$ws.Range("E29").Value = "  +7.20%  "
$ws.Range("E30").Value = "  +0.14%  "
$c = $ws.Range("D31")
$c.Value = "'3.96"
$c.ClearFormats()
$ws.Range("E31").Value = "  +0.62%  "
$c = $ws.Range("D32")
$c.Value = "'0.0547"
$c.ClearFormats()
$ws.Range("E32").Value = "  +0.76%  "
$c = $ws.Range("D33")
$c.Value = "'3.90"
$c.ClearFormats()
$ws.Range("E33").Value = "  -3.69%  "
$ws.Range("E34").Value = "  -7.81%  "
$ws.Range("E35").Value = "  +3.81%  "
$c = $ws.Range("D36")
$c.Value = "'0.675"
$c.ClearFormats()
$ws.Range("E36").Value = "  -1.00%  "
$c = $ws.Range("D37")
$c.Value = "'90.23"
$c.ClearFormats()
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "1.303.60"
$ws.Range("E39").Value = "  -2.91%  "
$c = $ws.Range("D40")
$c.Value = "'0.0190"
$c.ClearFormats()
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("E41").Value = "  -0.84%  "
$c = $ws.Range("D42")
$c.Value = "'0.955"
$c.ClearFormats()
$ws.Range("E42").Value = "  -6.77%  "
$ws.Range("E43").Value = "  -5.45%  "
$c = $ws.Range("D44")
$c.Value = "'2.73"
$c.ClearFormats()
$ws.Range("E44").Value = "  -3.34%  "
$c = $ws.Range("D45")
$c.Value = "'2.19"
$c.ClearFormats()
$ws.Range("E45").Value = "  -13.74%  "
$c = $ws.Range("D46")
$c.Value = "'6.12"
$c.ClearFormats()
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("D48").Value = "1.983.31"
$ws.Range("E48").Value = "  -2.15%  "
$c = $ws.Range("D50")
$c.Value = "'0.0667"
$c.ClearFormats()
$ws.Range("E50").Value = "  +7.24%  "
$c = $ws.Range("D51")
$c.Value = "'98.58"
$c.ClearFormats()
$ws.Range("E51").Value = "  -5.71%  "
